$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing data row (row 2) to the new row 3 so
# number formats / alignment / styles match (date format on A, 2-decimal
# number format on B, wrap-text style on C).
$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("A3:C3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New timesheet entry values
$ws.Range("A3").Value = [DateTime]"2014-11-05"
$ws.Range("B3").Value = 3.25
$ws.Range("C3").Value = "Setup css style sheet structures, start building new page layout in HTML - incorporating Boiler plate HTML and integrating XILIR layouts. Clean up HTML markup."

# The new row has taller (wrapped) text, so give it an explicit row height
$ws.Rows.Item(3).RowHeight = 30

# Update the active selection to reflect the new last cell, as in the diff
$ws.Range("C4").Select() | Out-Null
